$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I0/IF headers in I1/J1, matching the existing header style (s="1") ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Row 2 is a special case: I2 = 8, J2 = 9 ---
$ws.Cells.Item(2, 9).Value2 = 8
$ws.Cells.Item(2, 10).Value2 = 9

# --- Rows 3-33: I = 1, J = copy of H value for that row ---
for ($r = 3; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
